# AFDP-4158 Document Level - ACL Implementation
# File/Folder Assignment Rules
#
# Inserts a new "File - default access" rule row directly above the existing
# "Folder - default access" row in the Assignment Rules table (Sheet1),
# pushing every row below it down by one, and moves the on-screen selection
# to reflect where the author was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Folder - default access" rule currently lives on row 25 - insert a
# fresh blank row above it; Excel shifts rows 25:35 down to 26:36 and carries
# each row's own formatting with it.
$ws.Rows("25:25").Insert()

# Match the new row's formatting to the (now shifted-down) "Folder - default
# access" row directly below it, since both rows share the same layout.
$ws.Range("A26:H26").Copy()
$ws.Range("A25:H25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new "File - default access" rule.
$ws.Range("B25").Value = "File - default access"
$ws.Range("C25").Value = "FILE"
$ws.Range("D25").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G25").Value = "*, *"

# Reflect the author's final scroll position / selection on the sheet.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
